$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update test run results with the latest scenario run
$ws.Range("C2").Value = "Passed"
$ws.Range("H2").Value = "05_05_2020--23_36_02 897"

$ws.Range("B4").Value = "Yes"
$ws.Range("H4").Value = "05_05_2020--23_36_56 694"

$ws.Range("H6").Value = "05_05_2020--23_37_26 792"

# Update selection to reflect current active cell
$ws.Range("B5").Select()
